# kamlesh: updating URL method to extract based on keyword
# Adds a new "URL" worksheet (after "Login") holding a keyword -> value
# lookup table used by the URL-extraction logic, and makes it the active tab.

$wb = $excel.ActiveWorkbook

# --- Login sheet: it stops being the active/selected tab, and its last
#     selection moves to B1 (rather than B2). ---
$loginSheet = $wb.Worksheets.Item("Login")
$loginSheet.Range("B1").Select() | Out-Null

# --- Add the new "URL" worksheet right after "Login" ---
$newSheet = $wb.Worksheets.Add($null, $loginSheet)
$newSheet.Name = "URL"

# --- Header row + keyword/value rows ---
$newSheet.Range("A1").Value = "URL"
$newSheet.Range("B1").Value = "Value"
$newSheet.Range("A2").Value = "amazon"
$newSheet.Range("B2").Value = "https://amazon.co.in"
$newSheet.Range("A3").Value = "google"
$newSheet.Range("B3").Value = "google.co.in"

# --- B2 carries a real hyperlink (to the amazon URL) with the built-in
#     "Hyperlink" cell style (blue/underlined), like the Login sheet's
#     mailto link already does for A2. ---
$newSheet.Hyperlinks.Add($newSheet.Range("B2"), "https://amazon.co.in") | Out-Null
$newSheet.Range("B2").Style = "Hyperlink"

# --- Column B is widened to fit its longest value ("https://amazon.co.in"). ---
$newSheet.Columns.Item(2).ColumnWidth = 13.6

# --- Final selection on the new sheet sits on B3, and the URL sheet is the
#     active/selected tab. ---
$newSheet.Range("B3").Select() | Out-Null
